$d = $word.ActiveDocument

# WdHeaderFooterIndex constants
$wdHeaderFooterPrimary = 1
$wdHeaderFooterFirstPage = 2
$wdReplaceOne = 1

foreach ($sec in $d.Sections) {

    # --- Footers: "...2013-2018, FIX Protocol, Limited" -> "...2013-2020..." ---
    # Narrow match on just the two-digit year tail so we disturb as little of the
    # surrounding run structure as possible.
    foreach ($idx in @($wdHeaderFooterPrimary, $wdHeaderFooterFirstPage)) {
        $footer = $sec.Footers.Item($idx)
        if ($footer.Exists) {
            $footer.Range.Find.Execute("-2018,", $true, $false, $false, $false, $false,
                                        $true, 1, $false, "-2020,", $wdReplaceOne)
        }
    }

    # --- Primary header: "July 2018" -> "November 2020" ---
    $header = $sec.Headers.Item($wdHeaderFooterPrimary)
    if ($header.Exists) {
        $header.Range.Find.Execute("July 2018", $true, $false, $false, $false, $false,
                                    $true, 1, $false, "November 2020", $wdReplaceOne)
    }
}
